# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
# All cells in this workbook are plain text (dates/times/percentages are
# stored as literal strings, not real Excel dates/numbers), so every new
# cell's NumberFormat is forced to "@" (Text) before the value is written —
# this stops Excel's COM layer from auto-coercing strings like
# "2026-02-01" or "77.4%" into date/number values.

function Add-LogRows {
    param($Sheet, $StartRow, $Rows)

    $endRow = $StartRow + $Rows.Count - 1

    # Force every new cell to Text format first so values like
    # "2026-02-01" and "77.4%" are kept as literal strings.
    $Sheet.Range($Sheet.Cells.Item($StartRow, 1), $Sheet.Cells.Item($endRow, 6)).NumberFormat = "@"

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $row = $Rows[$i]
        for ($c = 1; $c -le 6; $c++) {
            $Sheet.Cells.Item($r, $c).Value = $row[$c - 1]
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: append rows 56-68
# ---------------------------------------------------------------------
$pirSheet = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-01","19:57:04","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:07","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:07","19:00","Bathroom","Motion Detected","Active"),
    @("2026-02-01","19:57:15","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:20","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:25","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:30","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:35","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:40","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:45","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:50","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:57:55","19:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","19:58:00","19:00","Bathroom","No Motion","Inactive")
)
Add-LogRows $pirSheet 56 $pirRows

# ---------------------------------------------------------------------
# Humidity sheet: append rows 45-55
# ---------------------------------------------------------------------
$humiditySheet = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-01","19:57:03","19:00","Bathroom","77.4%","Active"),
    @("2026-02-01","19:57:06","19:00","Bathroom","78.5%","Active"),
    @("2026-02-01","19:57:12","19:00","Bathroom","77.5%","Active"),
    @("2026-02-01","19:57:16","19:00","Bathroom","78.0%","Active"),
    @("2026-02-01","19:57:21","19:00","Bathroom","77.0%","Active"),
    @("2026-02-01","19:57:36","19:00","Bathroom","78.2%","Active"),
    @("2026-02-01","19:57:41","19:00","Bathroom","77.5%","Active"),
    @("2026-02-01","19:57:46","19:00","Bathroom","78.5%","Active"),
    @("2026-02-01","19:57:51","19:00","Bathroom","77.8%","Active"),
    @("2026-02-01","19:57:57","19:00","Bathroom","78.3%","Active"),
    @("2026-02-01","19:58:01","19:00","Bathroom","77.2%","Active")
)
Add-LogRows $humiditySheet 45 $humidityRows

# ---------------------------------------------------------------------
# Temperature sheet: append rows 45-55
# ---------------------------------------------------------------------
$temperatureSheet = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-01","19:57:03","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:06","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:12","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:16","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:21","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:37","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:42","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:47","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:52","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:57:57","19:00","Bathroom","25.2C","Active"),
    @("2026-02-01","19:58:02","19:00","Bathroom","25.2C","Active")
)
Add-LogRows $temperatureSheet 45 $temperatureRows

Write-Host "Appended $($pirRows.Count) rows to PIR, $($humidityRows.Count) rows to Humidity, $($temperatureRows.Count) rows to Temperature."
